$d = $word.ActiveDocument

# Re-type the single-word cells to trigger Word's proofing engine to drop
# the spellStart/spellEnd proofErr bookmarks that wrapped them.
$words = @("Mergesort", "ExactMatch", "FibIte", "MaxMinRec", "MultMatrixCost", "StackArray")
foreach ($w in $words) {
    $rng = $d.Content
    $rng.Find.Execute($w, $true, $true, $false, $false, $false, $true, 1, $false, $w, 2) | Out-Null
}

# Merge the "2" and the trailing spaces run into a single run.
$d.Content.Find.Execute("2       ", $true, $false, $false, $false, $false, $true, 1, $false, "2       ", 2) | Out-Null
